$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("documentos")

$ws.Range("B2").Value = "Documento contable nº. 1"
$ws.Range("C2").Value = "usuario1"
$ws.Range("D2").Value = "Tue Apr 09 13:50:25 CEST 2019"
$ws.Range("E2").Value = "DOCUMENTO_CONTABLE"
